# 개발일정관리.xlsx - development schedule updates
# Items #7 ("Service / 입력 값 전송 처리"), #8 ("XML / 값 입력 SQL") slipped
# from 2022-07-08 to 2022-07-11, and item #9 ("Login / 로그인 화면 - Front")
# now has its start/end dates filled in and is marked complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (item #7): start/end date moved 2022-07-08 -> 2022-07-11
$ws.Range("G11").Value = 44753
$ws.Range("H11").Value = 44753

# Row 12 (item #8): start/end date moved 2022-07-08 -> 2022-07-11
$ws.Range("G12").Value = 44753
$ws.Range("H12").Value = 44753

# Row 13 (item #9): start/end date now set to 2022-07-11 and status -> complete (Y)
$ws.Range("G13").Value = 44753
$ws.Range("H13").Value = 44753
$ws.Range("I13").Value = "Y"

# Best-effort: scroll the saved view down so row 7 is at the top of the window.
try {
    $excel.ActiveWindow.ScrollRow = 7
} catch {
}
